$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix beam motion systematic error values (row 27, columns C and D)
$ws.Range("C27").Value = 115
$ws.Range("D27").Value = 115

# Update the selected cell/active cell to reflect where the author left off
$ws.Range("G26").Select()
